# Apply the "updated merged file, version number, and queries" edit.
#
# Underlying fact change: Costco's count for "allergy question" went
# from 4 to 5. That single fact is duplicated (by the original SPARQL
# merge/combine process) across several sheets, so each copy - and any
# dependent SUM formula - needs to be bumped in step:
#   - Costco!B2            4 -> 5
#   - 'SPARQL query result'!C2   4 -> 5
#   - combine!C2            4 -> 5   (+ combine!C18 SUM 16 -> 17)
#   - 'combine color coded'!C3   4 -> 5   (+ E3 SUM 10 -> 11, C19 SUM 16 -> 17)
#
# On top of the data fix, the author's selection/active-sheet state when
# they saved the workbook moved from the "combine" tab to the
# "combine color coded" tab, with new active cells on several sheets.

$wb = $excel.ActiveWorkbook

# --- Costco: the root data point -----------------------------------------
$wsCostco = $wb.Worksheets.Item("Costco")
$wsCostco.Range("B2").Value = 5
$wsCostco.Activate()
$wsCostco.Range("B2").Select() | Out-Null

# --- SPARQL query result: mirrors the Costco count in column C -----------
$wsResult = $wb.Worksheets.Item("SPARQL query result")
$wsResult.Range("C2").Value = 5
$wsResult.Activate()
$wsResult.Range("C2").Select() | Out-Null

# --- combine: mirrors the Costco count in column C, recompute totals -----
$wsCombine = $wb.Worksheets.Item("combine")
$wsCombine.Range("C2").Value = 5
$wsCombine.Activate()
$wsCombine.Range("C2").Select() | Out-Null

# --- combine color coded: mirrors the Costco count in column C, recompute
#     totals; this is the sheet left active/selected when the file was saved
$wsColor = $wb.Worksheets.Item("combine color coded")
$wsColor.Range("C3").Value = 5
$wsColor.Activate()
$wsColor.Range("I10").Select() | Out-Null
